$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case rows, mirroring the structure of the existing row 2
# (Symbol, From_Time, To_Time, Interval, Trade_Amount, Take_Profit_PCT,
#  Stop_Loss_PCT, Fees_PCT), each with its own Test_Num and Symbol.

# Excel serial date values (whole days, no time-of-day component):
# 44378 = 2021-07-01, 44561 = 2021-12-31
$fromTime = 44378
$toTime   = 44561

$rows = @(
    @{ Row = 3; Num = 2; Symbol = "ETHUSD" },
    @{ Row = 4; Num = 3; Symbol = "DOTUSDT" },
    @{ Row = 5; Num = 4; Symbol = "BTCUSD" },
    @{ Row = 6; Num = 5; Symbol = "ADAUSDT" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Num
    $ws.Cells.Item($row, 2).Value = $r.Symbol
    $ws.Cells.Item($row, 3).Value = $fromTime
    $ws.Cells.Item($row, 4).Value = $toTime
    $ws.Cells.Item($row, 5).Value = 30
    $ws.Cells.Item($row, 6).Value = 10000
    $ws.Cells.Item($row, 7).Value = 1.5
    $ws.Cells.Item($row, 8).Value = 1
    $ws.Cells.Item($row, 9).Value = 0.075
}

$ws.Range("D13").Select()
